$wb = $excel.ActiveWorkbook

# --- Yearly sheets (2006-2010): clear superseded county-suppressed rows and drop trailing state-level rows that no longer belong ---
$ws = $wb.Worksheets.Item("2006")
$ws.Range("C3").Value2 = $null
$ws.Range("D3").Value2 = $null
$ws.Range("E3").Value2 = $null
$ws.Range("C13").Value2 = $null
$ws.Range("D13").Value2 = $null
$ws.Range("E13").Value2 = $null
$ws.Rows("53:54").Delete()

$ws = $wb.Worksheets.Item("2007")
$ws.Range("C3").Value2 = $null
$ws.Range("D3").Value2 = $null
$ws.Range("E3").Value2 = $null
$ws.Range("C13").Value2 = $null
$ws.Range("D13").Value2 = $null
$ws.Range("E13").Value2 = $null
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2008")
$ws.Range("C13").Value2 = $null
$ws.Range("D13").Value2 = $null
$ws.Range("E13").Value2 = $null
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2009")
$ws.Range("C13").Value2 = $null
$ws.Range("D13").Value2 = $null
$ws.Range("E13").Value2 = $null
$ws.Rows("53:55").Delete()

$ws = $wb.Worksheets.Item("2010")
$ws.Range("C13").Value2 = $null
$ws.Range("D13").Value2 = $null
$ws.Range("E13").Value2 = $null
$ws.Rows("53:55").Delete()

# --- Aggregate sheet: re-computed incidence/prevalence table (row 3 "Alaska" dropped, all subsequent rows shift up and are recomputed) ---
$ws = $wb.Worksheets.Item("Aggregate")

$ws.Range("D2").Value2 = 1091026.5426309681
$ws.Range("E2").Value2 = 14.431363345952535
$ws.Range("A3").Value2 = 4.0
$ws.Range("B3").Value2 = "Arizona"
$ws.Range("C3").Value2 = 866487.1122159892
$ws.Range("D3").Value2 = 6631643.360124848
$ws.Range("E3").Value2 = 13.065948591657628
$ws.Range("A4").Value2 = 6.0
$ws.Range("B4").Value2 = "California"
$ws.Range("C4").Value2 = 4513981.265956005
$ws.Range("D4").Value2 = 36944762.16889343
$ws.Range("E4").Value2 = 12.218190078800035
$ws.Range("A5").Value2 = 9.0
$ws.Range("B5").Value2 = "Connecticut"
$ws.Range("C5").Value2 = 515371.889014145
$ws.Range("D5").Value2 = 3216911.134948612
$ws.Range("E5").Value2 = 16.02070642906888
$ws.Range("A6").Value2 = 10.0
$ws.Range("B6").Value2 = "Delaware"
$ws.Range("C6").Value2 = 37621.98624757164
$ws.Range("D6").Value2 = 206483.30947034978
$ws.Range("E6").Value2 = 18.220352213491626
$ws.Range("A7").Value2 = 11.0
$ws.Range("B7").Value2 = "District of Columbia"
$ws.Range("C7").Value2 = 109403.39658694674
$ws.Range("D7").Value2 = 550985.4432505644
$ws.Range("E7").Value2 = 19.855950447895008
$ws.Range("A8").Value2 = 13.0
$ws.Range("B8").Value2 = "Geogia"
$ws.Range("C8").Value2 = 1847944.0559254412
$ws.Range("D8").Value2 = 12211231.561340636
$ws.Range("E8").Value2 = 15.133150547859733
$ws.Range("A9").Value2 = 16.0
$ws.Range("B9").Value2 = "Idaho"
$ws.Range("C9").Value2 = 104996.84186401556
$ws.Range("D9").Value2 = 1168919.0900528755
$ws.Range("E9").Value2 = 8.982387468688366
$ws.Range("A10").Value2 = 17.0
$ws.Range("B10").Value2 = "Illinois"
$ws.Range("C10").Value2 = 1580895.706567785
$ws.Range("D10").Value2 = 12758370.67075068
$ws.Range("E10").Value2 = 12.391047002514844
$ws.Range("A11").Value2 = 18.0
$ws.Range("B11").Value2 = "Indiana"
$ws.Range("C11").Value2 = 1006366.2885382811
$ws.Range("D11").Value2 = 7837909.728326183
$ws.Range("E11").Value2 = 12.839727981317218
$ws.Range("A12").Value2 = 19.0
$ws.Range("B12").Value2 = "Iowa"
$ws.Range("C12").Value2 = 287608.6462125769
$ws.Range("D12").Value2 = 3410013.627202459
$ws.Range("E12").Value2 = 8.434237444632387
$ws.Range("A13").Value2 = 20.0
$ws.Range("B13").Value2 = "Kansas"
$ws.Range("C13").Value2 = 396146.7413020134
$ws.Range("D13").Value2 = 3428397.615133414
$ws.Range("E13").Value2 = 11.554865735332672
$ws.Range("A14").Value2 = 21.0
$ws.Range("B14").Value2 = "Kentucky"
$ws.Range("C14").Value2 = 413087.5759674031
$ws.Range("D14").Value2 = 2957081.5267942287
$ws.Range("E14").Value2 = 13.969434803349209
$ws.Range("A15").Value2 = 22.0
$ws.Range("B15").Value2 = "Louisiana"
$ws.Range("C15").Value2 = 282442.7614159923
$ws.Range("D15").Value2 = 2178496.230470162
$ws.Range("E15").Value2 = 12.96503328605878
$ws.Range("A16").Value2 = 23.0
$ws.Range("B16").Value2 = "Maine"
$ws.Range("C16").Value2 = 109119.52465686374
$ws.Range("D16").Value2 = 825220.5687211744
$ws.Range("E16").Value2 = 13.223073780863679
$ws.Range("A17").Value2 = 24.0
$ws.Range("B17").Value2 = "Maryland"
$ws.Range("C17").Value2 = 1002975.984780239
$ws.Range("D17").Value2 = 6754689.3152769
$ws.Range("E17").Value2 = 14.848587965575783
$ws.Range("A18").Value2 = 26.0
$ws.Range("B18").Value2 = "Michigan"
$ws.Range("C18").Value2 = 1625547.03354681
$ws.Range("D18").Value2 = 11990510.426796142
$ws.Range("E18").Value2 = 13.55694608224577
$ws.Range("A19").Value2 = 27.0
$ws.Range("B19").Value2 = "Minnesota"
$ws.Range("C19").Value2 = 115805.68122609484
$ws.Range("D19").Value2 = 1222301.5589915775
$ws.Range("E19").Value2 = 9.474395281115143
$ws.Range("A20").Value2 = 28.0
$ws.Range("B20").Value2 = "Mississippi"
$ws.Range("C20").Value2 = 515551.3439893284
$ws.Range("D20").Value2 = 3630502.6189303007
$ws.Range("E20").Value2 = 14.20055011945513
$ws.Range("A21").Value2 = 29.0
$ws.Range("B21").Value2 = "Missouri"
$ws.Range("C21").Value2 = 766863.5236230671
$ws.Range("D21").Value2 = 5518463.724450555
$ws.Range("E21").Value2 = 13.896322634600988
$ws.Range("A22").Value2 = 30.0
$ws.Range("B22").Value2 = "Montana"
$ws.Range("C22").Value2 = 102944.19611154444
$ws.Range("D22").Value2 = 1059004.3401108915
$ws.Range("E22").Value2 = 9.720847423605917
$ws.Range("A23").Value2 = 31.0
$ws.Range("B23").Value2 = "Nebraska"
$ws.Range("C23").Value2 = 204904.80981110327
$ws.Range("D23").Value2 = 2201247.7197266174
$ws.Range("E23").Value2 = 9.308575676187472
$ws.Range("A24").Value2 = 32.0
$ws.Range("B24").Value2 = "Nevada"
$ws.Range("C24").Value2 = 351541.34141214454
$ws.Range("D24").Value2 = 3232267.2275993405
$ws.Range("E24").Value2 = 10.875998692510342
$ws.Range("A25").Value2 = 33.0
$ws.Range("B25").Value2 = "New Hampshire"
$ws.Range("C25").Value2 = 107549.14772854405
$ws.Range("D25").Value2 = 886427.4615532869
$ws.Range("E25").Value2 = 12.132876337121333
$ws.Range("A26").Value2 = 34.0
$ws.Range("B26").Value2 = "New Jersey"
$ws.Range("C26").Value2 = 1168380.2370796714
$ws.Range("D26").Value2 = 8196056.390538524
$ws.Range("E26").Value2 = 14.255395295089016
$ws.Range("A27").Value2 = 35.0
$ws.Range("B27").Value2 = "New Mexico"
$ws.Range("C27").Value2 = 179120.75376638884
$ws.Range("D27").Value2 = 1497759.9561924192
$ws.Range("E27").Value2 = 11.959243070014148
$ws.Range("A28").Value2 = 36.0
$ws.Range("B28").Value2 = "New York"
$ws.Range("C28").Value2 = 2781426.375356151
$ws.Range("D28").Value2 = 17587680.554876994
$ws.Range("E28").Value2 = 15.814628692382476
$ws.Range("A29").Value2 = 38.0
$ws.Range("B29").Value2 = "North Dakota"
$ws.Range("C29").Value2 = 37699.56334368641
$ws.Range("D29").Value2 = 425778.2815509131
$ws.Range("E29").Value2 = 8.854271102406717
$ws.Range("A30").Value2 = 39.0
$ws.Range("B30").Value2 = "Ohio"
$ws.Range("C30").Value2 = 997197.7980750023
$ws.Range("D30").Value2 = 8133869.580413682
$ws.Range("E30").Value2 = 12.259820350159654
$ws.Range("A31").Value2 = 40.0
$ws.Range("B31").Value2 = "Oklahoma"
$ws.Range("C31").Value2 = 487287.2730343682
$ws.Range("D31").Value2 = 3491912.5849495833
$ws.Range("E31").Value2 = 13.954738590382087
$ws.Range("A32").Value2 = 41.0
$ws.Range("B32").Value2 = "Oregon"
$ws.Range("C32").Value2 = 281481.12543031445
$ws.Range("D32").Value2 = 2525767.163670863
$ws.Range("E32").Value2 = 11.14438137762546
$ws.Range("A33").Value2 = 42.0
$ws.Range("B33").Value2 = "Pennsylvania"
$ws.Range("C33").Value2 = 1905109.2449985659
$ws.Range("D33").Value2 = 13667687.409110986
$ws.Range("E33").Value2 = 13.938782677518697
$ws.Range("A34").Value2 = 44.0
$ws.Range("B34").Value2 = "Rhode Island"
$ws.Range("C34").Value2 = 145915.22512547293
$ws.Range("D34").Value2 = 907042.8718127855
$ws.Range("E34").Value2 = 16.086916027888698
$ws.Range("A35").Value2 = 48.0
$ws.Range("B35").Value2 = "Texas"
$ws.Range("C35").Value2 = 3420043.6562700826
$ws.Range("D35").Value2 = 26030067.579157054
$ws.Range("E35").Value2 = 13.13881973556073
$ws.Range("A36").Value2 = 49.0
$ws.Range("B36").Value2 = "Utah"
$ws.Range("C36").Value2 = 404738.42055012786
$ws.Range("D36").Value2 = 3963226.5196104664
$ws.Range("E36").Value2 = 10.21234639371328
$ws.Range("A37").Value2 = 50.0
$ws.Range("B37").Value2 = "Vermont"
$ws.Range("C37").Value2 = 89457.40950244437
$ws.Range("D37").Value2 = 646238.5548435429
$ws.Range("E37").Value2 = 13.842784345186956
$ws.Range("A38").Value2 = 51.0
$ws.Range("B38").Value2 = "Virginia"
$ws.Range("C38").Value2 = 979105.1491755865
$ws.Range("D38").Value2 = 7200938.098480222
$ws.Range("E38").Value2 = 13.59691106610442
$ws.Range("A39").Value2 = 53.0
$ws.Range("B39").Value2 = "Washington"
$ws.Range("C39").Value2 = 330137.8695229291
$ws.Range("D39").Value2 = 3063863.432174641
$ws.Range("E39").Value2 = 10.775214915131084
$ws.Range("A40").Value2 = 54.0
$ws.Range("B40").Value2 = "West Virginia"
$ws.Range("C40").Value2 = 238133.07028626621
$ws.Range("D40").Value2 = 1877224.0092325318
$ws.Range("E40").Value2 = 12.685383796237643
$ws.Range("A41").Value2 = 55.0
$ws.Range("B41").Value2 = "Wisconsin"
$ws.Range("C41").Value2 = 410614.6794413228
$ws.Range("D41").Value2 = 3890402.9564449224
$ws.Range("E41").Value2 = 10.554553963647647
$ws.Range("A42").Value2 = 56.0
$ws.Range("B42").Value2 = "Wyoming"
$ws.Range("C42").Value2 = 34970.686195598726
$ws.Range("D42").Value2 = 367927.95107951324
$ws.Range("E42").Value2 = 9.504764748911718

# Remove the two rows that no longer fit after the re-computation (old Wisconsin/Wyoming rows, now redundant)
$ws.Rows("43:44").Delete()

Write-Host "done"